$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 409, shifting existing rows 409-451 down to 410-452
$ws.Rows.Item(409).Insert()

# Populate the newly inserted row 409 with the new data point
$ws.Cells.Item(409, 1).Value = 4
$ws.Cells.Item(409, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(409, 3).Value = "Los Lagos"
$ws.Cells.Item(409, 4).Value = 45212
$ws.Cells.Item(409, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(409, 5).Value = 10
$ws.Cells.Item(409, 6).Value = "Fruta"
$ws.Cells.Item(409, 7).Value = 100108
$ws.Cells.Item(409, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(409, 9).Value = 100108002
$ws.Cells.Item(409, 10).Value = "Mango"
$ws.Cells.Item(409, 11).Value = "Sin especificar"
$ws.Cells.Item(409, 12).Value = "Primera"
$ws.Cells.Item(409, 13).Value = 100
$ws.Cells.Item(409, 14).Value = 13000
$ws.Cells.Item(409, 15).Value = 13000
$ws.Cells.Item(409, 16).Value = 13000
$ws.Cells.Item(409, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(409, 18).Value = "Brasil"
$ws.Cells.Item(409, 19).Value = 3250
$ws.Cells.Item(409, 20).Value = 4
